$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B9").Value = "bc9b8f96e23a23185b98ef25692e5259"
$ws.Range("B17").Value = "6d85147a839ea3ac55b5d6f004e7aef0"
$ws.Range("B29").Value = "932c1bb69bdb11eb9f071ba98aeb37d5"
$ws.Range("B126").Value = "0e8e40d03f4ab904969ee58f99895996"
$ws.Range("B136").Value = "be6e2d30939511ef90af7ddc245aa206"
$ws.Range("B160").Value = "f621395159f113b00bc595d954de725b"
$ws.Range("B170").Value = "f5569e21c019eabf61d4098eaeb84c2f"
$ws.Range("B176").Value = "a218a17fec5319ca2fa92f87d646c4a5"
$ws.Range("B184").Value = "d2bca64db70b45f84c1e58f4fd552c94"
$ws.Range("B201").Value = "2531d54e1305d7e741c22516d8b18a53"
$ws.Range("B229").Value = "91dfa6a0e7de8a19c1e6b5f4d5a5077c"
$ws.Range("B230").Value = "ebbee841ca6a678c084a2e3151beca24"
$ws.Range("B234").Value = "fc9299dca116d10869fb7d133851c3b6"
$ws.Range("B287").Value = "63a3c99d70478e877f10bc650e02f22e"
$ws.Range("B299").Value = "c84bb76aa5ad595b32977ac03ac6772c"
$ws.Range("B308").Value = "fda0c60c095d94ab156b61dcb5489d43"
$ws.Range("B345").Value = "02d3049ffdaefb2d544cfbd86b8790f8"
$ws.Range("B470").Value = "e11a9cc46bed1da5741005a0ef219f46"
$ws.Range("B489").Value = "e97c1f07d1f73bb3afb061b17b5f515e"
$ws.Range("B514").Value = "54cc6f947e5c8aa2c1dc81ef5833b89f"
$ws.Range("B531").Value = "8d7b620123f834b0a2d5044f98cf1391"
$ws.Range("B563").Value = "05cb9a477ba28a3cc052a6500361e58c"
$ws.Range("B566").Value = "b81f17e92838431e2621315dfb396577"
$ws.Range("B579").Value = "ee1a38d9beeaab43b9f338874817fdb8"
$ws.Range("B588").Value = "9328d45f9fb021e84564eb92401c9b66"
$ws.Range("B632").Value = "9aeaa959732c5aa0194239ff31a74db5"
$ws.Range("B643").Value = "025e5e49f19857a1bff472dbfbd4e079"
$ws.Range("B645").Value = "ecf70847a71ebc8a4d8c99c66e167b40"
$ws.Range("B681").Value = "1564e15a977b19f6a79c7faf8a4de6e6"
$ws.Range("B696").Value = "37e3147df8fce8804dcb90d112343dac"
$ws.Range("B701").Value = "7484a6751fcbede002d4edc488191669"
$ws.Range("B716").Value = "54dba3e87c32ca8bdf54aee14a147675"
$ws.Range("B719").Value = "d19486f4fe6118e3a48dbee74a9f737a"
$ws.Range("B731").Value = "afa32301fe7c3e656e0a5169bb5010c2"
$ws.Range("B845").Value = "0244a7d27ded1427cbf7477266280769"
$ws.Range("B848").Value = "5878ebe966667d1f8ef44bddc0589681"
$ws.Range("B853").Value = "92782e709e79bd5318a2e10eab8b64d0"
$ws.Range("B880").Value = "d64c9a973b8b4aa4e508a366f2c61478"
